# edit.ps1
# Applies the changes described by the commit "adding averages and more checks"
# to the "Training Dashboard" (sheet1) and "Exam Dashboard" (sheet2) worksheets.

$wb = $excel.ActiveWorkbook

$wsTraining = $wb.Worksheets.Item("Training Dashboard")
$wsExam     = $wb.Worksheets.Item("Exam Dashboard")

# ---------------------------------------------------------------------------
# 1) Header / title formatting: bold white text on the title row and the
#    column-header row (these share a single font/style definition across the
#    whole workbook, so apply consistently on both sheets).
# ---------------------------------------------------------------------------
$white = 16777215  # RGB(255,255,255)

$wsTraining.Range("A1").Font.Color = $white
$wsTraining.Range("A1").Font.Size = 11
$wsTraining.Range("A2:K2").Font.Color = $white

$wsExam.Range("A1").Font.Color = $white
$wsExam.Range("A1").Font.Size = 11
$wsExam.Range("A2:G2").Font.Color = $white

# ---------------------------------------------------------------------------
# 2) Training Dashboard: "PERIOD TO EXPIRE" (col H) decreases by 8 days and
#    "LAST UPDATE" (col I) moves from 08-Sep-2025 to 16-Sep-2025, for every
#    data row (3-16).
# ---------------------------------------------------------------------------
$periodToExpire = @{
    3  = 401
    4  = 403
    5  = 699
    6  = 446
    7  = 399
    8  = 394
    9  = 633
    10 = 256
    11 = 315
    12 = 316
    13 = 318
    14 = 314
    15 = 335
    16 = 336
}

# Make sure the LAST UPDATE column keeps storing the date as literal text
# (matching the existing data) instead of being auto-converted to a date
# serial number.
$wsTraining.Range("I3:I16").NumberFormat = "@"

foreach ($row in 3..16) {
    $wsTraining.Cells.Item($row, 8).Value = $periodToExpire[$row]
    $wsTraining.Cells.Item($row, 9).Value = "16-Sep-2025"
}

# ---------------------------------------------------------------------------
# 3) Exam Dashboard: widen the COMMENTS column and update the comment text
#    for the existing exam rows.
# ---------------------------------------------------------------------------
$wsExam.Columns.Item(5).ColumnWidth = 14.166666666666666

$wsExam.Range("E3").Value = "date is valid"
$wsExam.Range("E4").Value = "date is valid"
